$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.878023624420166
$ws.Range("B1").Value = 0.8574672937393188
$ws.Range("C1").Value = 0.8891956806182861
$ws.Range("D1").Value = 1.107177734375
$ws.Range("E1").Value = 0.9926592111587524
